$d = $word.ActiveDocument

function Replace-WithBreaks($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $replacement = $find.Replacement
    $replacement.ClearFormatting()
    $result = $find.Execute(
        $findText,   # FindText
        $false,      # MatchCase
        $false,      # MatchWholeWord
        $false,      # MatchWildcards
        $false,      # MatchSoundsLike
        $false,      # MatchAllWordForms
        $true,       # Forward
        1,           # Wrap (wdFindContinue)
        $false,      # Format
        $replaceText,# ReplaceWith
        2            # Replace (wdReplaceAll)
    )
    if (-not $result) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) Objetivos (PT) - split "Gerais..." paragraph into two sentences separated by a blank line
Replace-WithBreaks "necessárias.Específicos" "necessárias.^l^lEspecíficos"

# 2) Objetivos (EN) - split "Overview..." paragraph the same way
Replace-WithBreaks "necessary.Specific" "necessary.^l^lSpecific"

# 3) Programa resumido (PT) - split "estruturais" / "como estereoquímica" with a single break
Replace-WithBreaks "estruturaiscomo estereoquímica" "estruturais^lcomo estereoquímica"

# 4) Programa (PT) - split the numbered list into separate lines
Replace-WithBreaks "orgânicos2.Alcanos" "orgânicos^l2.Alcanos"
Replace-WithBreaks "Radicalar. 3.Isomeria" "Radicalar. ^l3.Isomeria"
Replace-WithBreaks "Espaciais.4.Haletos" "Espaciais.^l4.Haletos"
Replace-WithBreaks "E1, E2. 5.Alcenos" "E1, E2. ^l5.Alcenos"
Replace-WithBreaks "cinético) 6. Fundamentos" "cinético) ^l6. Fundamentos"
Replace-WithBreaks "Fluorescencia 7.Compostos" "Fluorescencia ^l7.Compostos"
Replace-WithBreaks "Nucleofílica.8.Álcoois" "Nucleofílica.^l8.Álcoois"

# 5) Programa (EN) - split the numbered list into separate lines
Replace-WithBreaks "compounds2.Alkanes" "compounds^l2.Alkanes"
Replace-WithBreaks "reaction.3.Constitutional" "reaction.^l3.Constitutional"
Replace-WithBreaks "stereoisomers.4.Alkyl" "stereoisomers.^l4.Alkyl"
Replace-WithBreaks "E1, E2.5.Alkenes" "E1, E2.^l5.Alkenes"
Replace-WithBreaks "product).6 Background" "product).^l6 Background"
Replace-WithBreaks "techniques7.Aromatic" "techniques^l7.Aromatic"
Replace-WithBreaks "Substitution.8.Alcohols" "Substitution.^l8.Alcohols"

# 6) Avaliação - split "Método" text into two sentences separated by a blank line
Replace-WithBreaks "letivoAos alunos" "letivo^l^lAos alunos"
